$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting rows 44-68 down to 45-69
$ws.Rows("44").Insert()

# Populate the new row 44 with data (same template as the row that was previously
# at position 44, now shifted to 45, but with updated date/price values)
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 44460
$ws.Range("D44").NumberFormat = $ws.Range("D45").NumberFormat
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112043
$ws.Range("G44").Value = "Pepino ensalada"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 16000
$ws.Range("L44").Value = 17000
$ws.Range("M44").Value = 16500
$ws.Range("N44").Value = "$/caja 60 unidades"
$ws.Range("O44").Value = "Región de Arica y Parinacota"
$ws.Range("P44").Value = 275
$ws.Range("Q44").Value = 60
$ws.Range("R44").Value = "Hortaliza"
